$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 897.9048
$ws.Range("I43").Value = 733.1667
$ws.Range("J43").Value = 963.8
$ws.Range("K43").Value = 733.1667
$ws.Range("L43").Value = 963.8
$ws.Range("M43").Value = -664.1667
$ws.Range("N43").Value = -1101.8
$ws.Range("H137").Value = 1341.84
$ws.Range("I137").Value = 883.1429000000001
$ws.Range("K137").Value = 2649.4287
$ws.Range("M137").Value = -99.42870000000039
$ws.Range("H138").Value = 3599.7585
$ws.Range("J138").Value = 3254.7693
$ws.Range("L138").Value = 9764.3079
$ws.Range("N138").Value = -20044.3079

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2222.425
$ws.Range("I61").Value = 1429.6487
$ws.Range("K61").Value = 1429.6487
$ws.Range("M61").Value = -1217.6487
$ws.Range("H74").Value = 1860.8334
$ws.Range("I74").Value = 1762.7142
$ws.Range("K74").Value = 1762.7142
$ws.Range("M74").Value = -888.7141999999999
$ws.Range("H77").Value = 1860.8334
$ws.Range("I77").Value = 1762.7142
$ws.Range("K77").Value = 8813.571
$ws.Range("M77").Value = -4445.571
$ws.Range("H132").Value = 1518.2778
$ws.Range("I132").Value = 998.74194
$ws.Range("J132").Value = 4739.4
$ws.Range("K132").Value = 2996.22582
$ws.Range("L132").Value = 14218.2
$ws.Range("M132").Value = -466.2258200000001
$ws.Range("N132").Value = -19278.2
$ws.Range("H135").Value = 30429
$ws.Range("J135").Value = 30429
$ws.Range("L135").Value = 30429
$ws.Range("N135").Value = -40569
$ws.Range("H136").Value = 2222.425
$ws.Range("I136").Value = 1429.6487
$ws.Range("K136").Value = 4288.9461
$ws.Range("M136").Value = -1738.9461
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 226888.67
$ws.Range("I86").Value = 9999.666999999999
$ws.Range("J86").Value = 335333.16
$ws.Range("K86").Value = 9999.666999999999
$ws.Range("L86").Value = 335333.16
$ws.Range("M86").Value = -8876.666999999999
$ws.Range("N86").Value = -337579.16
$ws.Range("H89").Value = 226888.67
$ws.Range("I89").Value = 9999.666999999999
$ws.Range("J89").Value = 335333.16
$ws.Range("K89").Value = 49998.335
$ws.Range("L89").Value = 1676665.8
$ws.Range("M89").Value = -44382.335
$ws.Range("N89").Value = -1687897.8

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1721.5385
$ws.Range("I31").Value = 1553.4
$ws.Range("J31").Value = 2282
$ws.Range("K31").Value = 1553.4
$ws.Range("L31").Value = 2282
$ws.Range("M31").Value = -1258.4
$ws.Range("N31").Value = -2872
$ws.Range("H34").Value = 1721.5385
$ws.Range("I34").Value = 1553.4
$ws.Range("J34").Value = 2282
$ws.Range("K34").Value = 1553.4
$ws.Range("L34").Value = 2282
$ws.Range("M34").Value = -1351.4
$ws.Range("N34").Value = -2686

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5546307.5
$ws.Range("I4").Value = 2512587.5
$ws.Range("J4").Value = 10400260
$ws.Range("K4").Value = 7537762.5
$ws.Range("L4").Value = 31200780
$ws.Range("M4").Value = -7537650.5
$ws.Range("N4").Value = -31201004
$ws.Range("H64").Value = 2746.1428
$ws.Range("J64").Value = 3128.5
$ws.Range("L64").Value = 9385.5
$ws.Range("N64").Value = -9925.5
$ws.Range("H67").Value = 2746.1428
$ws.Range("J67").Value = 3128.5
$ws.Range("L67").Value = 9385.5
$ws.Range("N67").Value = -11257.5
$ws.Range("H76").Value = 3300
$ws.Range("J76").Value = 4200
$ws.Range("L76").Value = 12600
$ws.Range("N76").Value = -13366
$ws.Range("H79").Value = 3300
$ws.Range("J79").Value = 4200
$ws.Range("L79").Value = 12600
$ws.Range("N79").Value = -15252
$ws.Range("H131").Value = 11236.091
$ws.Range("J131").Value = 12117.873
$ws.Range("L131").Value = 36353.619
$ws.Range("N131").Value = -46433.619
$ws.Range("I141").Value = 4370.6
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 13111.8
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -7931.800000000001
$ws.Range("N141").ClearContents()

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2416.6667
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 3125
$ws.Range("K80").Value = 1000
$ws.Range("L80").Value = 3125
$ws.Range("M80").Value = -2
$ws.Range("N80").Value = -5121
$ws.Range("H83").Value = 2416.6667
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 3125
$ws.Range("K83").Value = 5000
$ws.Range("L83").Value = 15625
$ws.Range("M83").Value = -8
$ws.Range("N83").Value = -25609
$ws.Range("H132").Value = 1133877.8
$ws.Range("I132").Value = 1924896.6
$ws.Range("J132").Value = 3850.7856
$ws.Range("K132").Value = 5774689.800000001
$ws.Range("L132").Value = 11552.3568
$ws.Range("M132").Value = -5772159.800000001
$ws.Range("N132").Value = -16612.3568

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3687.375
$ws.Range("J22").Value = 2785.5715
$ws.Range("L22").Value = 2785.5715
$ws.Range("N22").Value = -3375.5715
$ws.Range("H27").Value = 3687.375
$ws.Range("J27").Value = 2785.5715
$ws.Range("L27").Value = 2785.5715
$ws.Range("N27").Value = -2999.5715
$ws.Range("H127").Value = 75000
$ws.Range("J127").Value = 75000
$ws.Range("L127").Value = 75000
$ws.Range("N127").Value = -84920
$ws.Range("H132").Value = 1441.5254
$ws.Range("I132").Value = 950.4286
$ws.Range("J132").Value = 2157.7083
$ws.Range("K132").Value = 2851.2858
$ws.Range("L132").Value = 6473.124899999999
$ws.Range("M132").Value = -321.2857999999997
$ws.Range("N132").Value = -11533.1249

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 60962.46
$ws.Range("I122").Value = 71829.45
$ws.Range("J122").Value = 1194
$ws.Range("K122").Value = 215488.35
$ws.Range("L122").Value = 3582
$ws.Range("M122").Value = -213038.35
$ws.Range("N122").Value = -8482
$ws.Range("H132").Value = 1454.9487
$ws.Range("I132").Value = 787
$ws.Range("J132").Value = 3681.4443
$ws.Range("K132").Value = 2361
$ws.Range("L132").Value = 11044.3329
$ws.Range("M132").Value = 169
$ws.Range("N132").Value = -16104.3329
$ws.Range("H136").Value = 22224602
$ws.Range("I136").Value = 30866282
$ws.Range("K136").Value = 92598846
$ws.Range("M136").Value = -92596296
